$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'31.151.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.81%  "

$ws.Range("D3").Value = "'1.990.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.73%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'0.7877"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +66.58%  "

$ws.Range("D6").Value = "'254.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.26%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'0.3494"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +21.12%  "

$ws.Range("D9").Value = "'27.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +26.25%  "

$ws.Range("D10").Value = "'0.06984"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.85%  "

$ws.Range("D11").Value = "'0.8430"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.07%  "

$ws.Range("D12").Value = "'0.08187"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.64%  "

$ws.Range("D13").Value = "'1.990.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.78%  "

$ws.Range("D14").Value = "'100.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.83%  "

$ws.Range("D15").Value = "'5.595"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.53%  "

$ws.Range("D16").Value = "'15.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +16.75%  "

$ws.Range("D17").Value = "'273.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.15%  "

$ws.Range("D18").Value = "'31.149.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").Value = "'5.874"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.36%  "

$ws.Range("D20").Value = "'0.000007898"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.82%  "

$ws.Range("D21").Value = "'2.255.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.26%  "

$ws.Range("D22").Value = "'1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").Value = "'7.059"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.15%  "

$ws.Range("D25").Value = "'10.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.94%  "

$ws.Range("D26").Value = "'0.1503"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +54.85%  "

$ws.Range("D27").Value = "'164.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("D28").Value = "'19.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.27%  "

$ws.Range("D29").Value = "'2.323"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +21.21%  "

$ws.Range("E30").Value = "  +6.05%  "

$ws.Range("D31").Value = "'1.362"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.48%  "

$ws.Range("D32").Value = "'4.583"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.31%  "

$ws.Range("D33").Value = "'4.412"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.08%  "

$ws.Range("D34").Value = "'0.05217"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.57%  "

$ws.Range("D35").Value = "'1.225"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.31%  "

$ws.Range("D36").Value = "'0.7786"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.71%  "

$ws.Range("D37").Value = "'2.759"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").Value = "'0.02004"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.59%  "

$ws.Range("D39").Value = "'2.896"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").Value = "'6.621"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.22%  "

$ws.Range("D41").Value = "'79.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.55%  "

$ws.Range("D42").Value = "'0.4657"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.22%  "

$ws.Range("D43").Value = "'2.118"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.70%  "

$ws.Range("D44").Value = "'105.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.48%  "

$ws.Range("D45").Value = "'0.8478"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.18%  "

$ws.Range("D46").Value = "'1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").Value = "'7.672"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.96%  "

$ws.Range("D48").Value = "'9.918"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.25%  "

$ws.Range("D49").Value = "'36.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.64%  "

$ws.Range("D50").Value = "'0.4297"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.51%  "

$ws.Range("D51").Value = "'1.530"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.21%  "

